# EPICP-1: changed DPE and DD for Tracy and Ines with regards to variable names,
# changed DPE for Franzi to remove some leftovers in the variable column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the input_variables (column F) and a few algorithm (column H) / comment (column I)
# cells to use the new variable naming convention.
$ws.Range("F2").Value  = "bmi0"
$ws.Range("F3").Value  = "GJ"
$ws.Range("H3").Value  = "GJ*4.184"
$ws.Range("I3").Value  = "GJ [kJ] calculated in kcal"

$ws.Range("F12").Value = "waist0"
$ws.Range("F13").Value = "hip0"

$ws.Range("F20").Value = "ZK"
$ws.Range("F21").Value = "ZE"
$ws.Range("F22").Value = "ZF"
$ws.Range("F23").Value = "ZA"
$ws.Range("F24").Value = "ZB"
$ws.Range("F25").Value = "FS"
$ws.Range("F26").Value = "FU"
$ws.Range("F27").Value = "FP"
$ws.Range("F28").Value = "KD;KM"
$ws.Range("H28").Value = "KD+KM"
$ws.Range("F31").Value = "KMT"
$ws.Range("F32").Value = "KMF"
$ws.Range("F35").Value = "MNA"
$ws.Range("F36").Value = "MNA;MK"
$ws.Range("H36").Value = "MNA/MK"

# Update the view: scroll/zoom position changed as part of the edit session
$window = $excel.ActiveWindow
$window.Zoom = 70
$window.ScrollColumn = 4
$ws.Range("E5").Select()
